$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"

# Zero out the billed amount / pricing figures (no-violation scenario)
$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
